$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.274.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.669.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.58%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.24%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.546"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.668.09"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.158"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.358"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.164.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("E16").Value = "  -2.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.231.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.675.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "362.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.04%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.802.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("E29").Value = "  -2.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.992"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "552.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "155.72"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.53%  "
$ws.Range("E40").Value = "  -2.56%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.99%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("E44").Value = "  -4.85%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₆0299"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.586"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "153.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.29%  "
